# Add a new "to buy" line item to the bottom of the list (row 15):
#   B15 = "Electric Bike Scooter Funn"
#   C15 = "?"
# and leave the selection parked on the newly-added price cell, matching
# the author's final cursor position/scroll state in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "Electric Bike Scooter Funn"
$ws.Range("C15").Value = "?"

# Selecting the cell also clears the stale topLeftCell/selection scroll
# state that was left over from the previous edit session.
$ws.Range("C15").Select()
